$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.731.88"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "2.943.59"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "2.942.12"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.01"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "3.430.21"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "62.765.05"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "2.945.05"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.77"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.14%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "0.0₃0997"
$ws.Range("E32").Value = "  +13.42%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.21"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.991"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.59"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.37"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.276"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.34"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("D46").Value = "2.680.13"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0337"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "354.98"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.59%  "
